$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.823.17"
$ws.Range("E2").Value = "  -5.25%  "
$ws.Range("D3").Value = "3.364.44"
$ws.Range("E3").Value = "  -6.59%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.595"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.29%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.354.52"
$ws.Range("E9").Value = "  -6.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -12.71%  "
$ws.Range("E11").Value = "  -7.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -11.35%  "
$ws.Range("E13").Value = "  -11.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.45%  "
$ws.Range("D15").Value = "3.901.44"
$ws.Range("E15").Value = "  -6.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "603.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.13%  "
$ws.Range("D17").Value = "66.848.20"
$ws.Range("E17").Value = "  -5.35%  "
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").Value = "3.368.74"
$ws.Range("E19").Value = "  -7.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.913"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.21%  "
$ws.Range("E26").Value = "  -9.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -14.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.71%  "
$ws.Range("E34").Value = "  -8.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "534.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.24%  "
$ws.Range("D37").Value = "3.762.08"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +42.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("D41").Value = "0.0₃0729"
$ws.Range("E41").Value = "  -14.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.353"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.127"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.28%  "
$ws.Range("E46").Value = "  -10.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.94%  "
$ws.Range("E48").Value = "  -12.68%  "
$ws.Range("E49").Value = "  -8.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.11%  "
